# "final pre-test-rig animation sketches"
# Trim the trailing repeated animation frames (rows 13-16) from the
# pulsar sketch sheet, leaving row 12 as the last frame, and move the
# active selection to where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove animation frame rows 13-16 (row 12 keeps its existing values;
# Excel reflows dimension/used-range automatically on row delete).
$ws.Rows("13:16").Delete()

# Restore the author's last selection/cursor position on the sheet.
$ws.Range("AR23").Select()
